# ---------------------------------------------------------------
# feat: add 2022-Q1 data
# ---------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# 1) Insert a new "2022-Q1" worksheet right before the "总计" sheet,
#    using the existing "2021-Q4" fund-holding sheet as a formatting
#    template (bold/centered/bordered header row + A-column index).
$total = $wb.Worksheets.Item("总计")
$ref = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($total)
$q1.Name = "2022-Q1"

$ref.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$ref.Range("A2:A16").Copy()
$q1.Range("A2:A32").PasteSpecial(-4122)

# Header row (same layout as the other quarterly fund-holding sheets)
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# 31 fund rows for 2022-Q1
$q1.Range("A2").Value = 0
$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "004854"
$q1.Range("C2").Value = "广发中证全指汽车指数A"
$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value = "22.01"
$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value = "94.43"
$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value = "6.99"
$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value = "1.5385"
$q1.Range("H2").Value = 5

$q1.Range("A3").Value = 1
$q1.Range("B3").NumberFormat = "@"
$q1.Range("B3").Value = "320006"
$q1.Range("C3").Value = "诺安灵活配置混合"
$q1.Range("D3").NumberFormat = "@"
$q1.Range("D3").Value = "10.15"
$q1.Range("E3").NumberFormat = "@"
$q1.Range("E3").Value = "70.73"
$q1.Range("F3").NumberFormat = "@"
$q1.Range("F3").Value = "7.53"
$q1.Range("G3").NumberFormat = "@"
$q1.Range("G3").Value = "0.7643"
$q1.Range("H3").Value = 1

$q1.Range("A4").Value = 2
$q1.Range("B4").NumberFormat = "@"
$q1.Range("B4").Value = "213001"
$q1.Range("C4").Value = "宝盈鸿利收益灵活配置混合A"
$q1.Range("D4").NumberFormat = "@"
$q1.Range("D4").Value = "17.98"
$q1.Range("E4").NumberFormat = "@"
$q1.Range("E4").Value = "90.37"
$q1.Range("F4").NumberFormat = "@"
$q1.Range("F4").Value = "4.12"
$q1.Range("G4").NumberFormat = "@"
$q1.Range("G4").Value = "0.7408"
$q1.Range("H4").Value = 6

$q1.Range("A5").Value = 3
$q1.Range("B5").NumberFormat = "@"
$q1.Range("B5").Value = "004855"
$q1.Range("C5").Value = "广发中证全指汽车指数C"
$q1.Range("D5").NumberFormat = "@"
$q1.Range("D5").Value = "6.11"
$q1.Range("E5").NumberFormat = "@"
$q1.Range("E5").Value = "94.43"
$q1.Range("F5").NumberFormat = "@"
$q1.Range("F5").Value = "6.99"
$q1.Range("G5").NumberFormat = "@"
$q1.Range("G5").Value = "0.4271"
$q1.Range("H5").Value = 5

$q1.Range("A6").Value = 4
$q1.Range("B6").NumberFormat = "@"
$q1.Range("B6").Value = "013610"
$q1.Range("C6").Value = "中信保诚前瞻优势混合"
$q1.Range("D6").NumberFormat = "@"
$q1.Range("D6").Value = "16.85"
$q1.Range("E6").NumberFormat = "@"
$q1.Range("E6").Value = "58.37"
$q1.Range("F6").NumberFormat = "@"
$q1.Range("F6").Value = "1.33"
$q1.Range("G6").NumberFormat = "@"
$q1.Range("G6").Value = "0.2241"
$q1.Range("H6").Value = 10

$q1.Range("A7").Value = 5
$q1.Range("B7").NumberFormat = "@"
$q1.Range("B7").Value = "001543"
$q1.Range("C7").Value = "宝盈新锐灵活配置混合A"
$q1.Range("D7").NumberFormat = "@"
$q1.Range("D7").Value = "3.21"
$q1.Range("E7").NumberFormat = "@"
$q1.Range("E7").Value = "93.26"
$q1.Range("F7").NumberFormat = "@"
$q1.Range("F7").Value = "4.89"
$q1.Range("G7").NumberFormat = "@"
$q1.Range("G7").Value = "0.1570"
$q1.Range("H7").Value = 5

$q1.Range("A8").Value = 6
$q1.Range("B8").NumberFormat = "@"
$q1.Range("B8").Value = "516110"
$q1.Range("C8").Value = "国泰中证800汽车与零部件ETF"
$q1.Range("D8").NumberFormat = "@"
$q1.Range("D8").Value = "1.74"
$q1.Range("E8").NumberFormat = "@"
$q1.Range("E8").Value = "97.85"
$q1.Range("F8").NumberFormat = "@"
$q1.Range("F8").Value = "4.04"
$q1.Range("G8").NumberFormat = "@"
$q1.Range("G8").Value = "0.0703"
$q1.Range("H8").Value = 8

$q1.Range("A9").Value = 7
$q1.Range("B9").NumberFormat = "@"
$q1.Range("B9").Value = "000664"
$q1.Range("C9").Value = "国联安通盈灵活配置混合A"
$q1.Range("D9").NumberFormat = "@"
$q1.Range("D9").Value = "7.87"
$q1.Range("E9").NumberFormat = "@"
$q1.Range("E9").Value = "24.12"
$q1.Range("F9").NumberFormat = "@"
$q1.Range("F9").Value = "0.71"
$q1.Range("G9").NumberFormat = "@"
$q1.Range("G9").Value = "0.0559"
$q1.Range("H9").Value = 8

$q1.Range("A10").Value = 8
$q1.Range("B10").NumberFormat = "@"
$q1.Range("B10").Value = "000538"
$q1.Range("C10").Value = "诺安优势行业灵活配置混合A"
$q1.Range("D10").NumberFormat = "@"
$q1.Range("D10").Value = "0.67"
$q1.Range("E10").NumberFormat = "@"
$q1.Range("E10").Value = "76.41"
$q1.Range("F10").NumberFormat = "@"
$q1.Range("F10").Value = "8.33"
$q1.Range("G10").NumberFormat = "@"
$q1.Range("G10").Value = "0.0558"
$q1.Range("H10").Value = 3

$q1.Range("A11").Value = 9
$q1.Range("B11").NumberFormat = "@"
$q1.Range("B11").Value = "001228"
$q1.Range("C11").Value = "国联安鑫享灵活配置混合A"
$q1.Range("D11").NumberFormat = "@"
$q1.Range("D11").Value = "5.71"
$q1.Range("E11").NumberFormat = "@"
$q1.Range("E11").Value = "22.58"
$q1.Range("F11").NumberFormat = "@"
$q1.Range("F11").Value = "0.65"
$q1.Range("G11").NumberFormat = "@"
$q1.Range("G11").Value = "0.0371"
$q1.Range("H11").Value = 7

$q1.Range("A12").Value = 10
$q1.Range("B12").NumberFormat = "@"
$q1.Range("B12").Value = "005571"
$q1.Range("C12").Value = "中银证券新能源灵活配置混合A"
$q1.Range("D12").NumberFormat = "@"
$q1.Range("D12").Value = "0.91"
$q1.Range("E12").NumberFormat = "@"
$q1.Range("E12").Value = "90.25"
$q1.Range("F12").NumberFormat = "@"
$q1.Range("F12").Value = "3.97"
$q1.Range("G12").NumberFormat = "@"
$q1.Range("G12").Value = "0.0361"
$q1.Range("H12").Value = 8

$q1.Range("A13").Value = 11
$q1.Range("B13").NumberFormat = "@"
$q1.Range("B13").Value = "519616"
$q1.Range("C13").Value = "银河君信灵活配置混合A"
$q1.Range("D13").NumberFormat = "@"
$q1.Range("D13").Value = "4.54"
$q1.Range("E13").NumberFormat = "@"
$q1.Range("E13").Value = "24.42"
$q1.Range("F13").NumberFormat = "@"
$q1.Range("F13").Value = "0.74"
$q1.Range("G13").NumberFormat = "@"
$q1.Range("G13").Value = "0.0336"
$q1.Range("H13").Value = 8

$q1.Range("A14").Value = 12
$q1.Range("B14").NumberFormat = "@"
$q1.Range("B14").Value = "519618"
$q1.Range("C14").Value = "银河君信灵活配置混合I"
$q1.Range("D14").NumberFormat = "@"
$q1.Range("D14").Value = "4.54"
$q1.Range("E14").NumberFormat = "@"
$q1.Range("E14").Value = "24.42"
$q1.Range("F14").NumberFormat = "@"
$q1.Range("F14").Value = "0.74"
$q1.Range("G14").NumberFormat = "@"
$q1.Range("G14").Value = "0.0336"
$q1.Range("H14").Value = 8

$q1.Range("A15").Value = 13
$q1.Range("B15").NumberFormat = "@"
$q1.Range("B15").Value = "002186"
$q1.Range("C15").Value = "国联安鑫享灵活配置混合C"
$q1.Range("D15").NumberFormat = "@"
$q1.Range("D15").Value = "5.14"
$q1.Range("E15").NumberFormat = "@"
$q1.Range("E15").Value = "22.58"
$q1.Range("F15").NumberFormat = "@"
$q1.Range("F15").Value = "0.65"
$q1.Range("G15").NumberFormat = "@"
$q1.Range("G15").Value = "0.0334"
$q1.Range("H15").Value = 7

$q1.Range("A16").Value = 14
$q1.Range("B16").NumberFormat = "@"
$q1.Range("B16").Value = "002485"
$q1.Range("C16").Value = "国联安通盈灵活配置混合C"
$q1.Range("D16").NumberFormat = "@"
$q1.Range("D16").Value = "4.54"
$q1.Range("E16").NumberFormat = "@"
$q1.Range("E16").Value = "24.12"
$q1.Range("F16").NumberFormat = "@"
$q1.Range("F16").Value = "0.71"
$q1.Range("G16").NumberFormat = "@"
$q1.Range("G16").Value = "0.0322"
$q1.Range("H16").Value = 8

$q1.Range("A17").Value = 15
$q1.Range("B17").NumberFormat = "@"
$q1.Range("B17").Value = "007581"
$q1.Range("C17").Value = "宝盈鸿利收益灵活配置混合C"
$q1.Range("D17").NumberFormat = "@"
$q1.Range("D17").Value = "0.73"
$q1.Range("E17").NumberFormat = "@"
$q1.Range("E17").Value = "90.37"
$q1.Range("F17").NumberFormat = "@"
$q1.Range("F17").Value = "4.12"
$q1.Range("G17").NumberFormat = "@"
$q1.Range("G17").Value = "0.0301"
$q1.Range("H17").Value = 6

$q1.Range("A18").Value = 16
$q1.Range("B18").NumberFormat = "@"
$q1.Range("B18").Value = "002053"
$q1.Range("C18").Value = "诺安优势行业灵活配置混合C"
$q1.Range("D18").NumberFormat = "@"
$q1.Range("D18").Value = "0.30"
$q1.Range("E18").NumberFormat = "@"
$q1.Range("E18").Value = "76.41"
$q1.Range("F18").NumberFormat = "@"
$q1.Range("F18").Value = "8.33"
$q1.Range("G18").NumberFormat = "@"
$q1.Range("G18").Value = "0.0250"
$q1.Range("H18").Value = 3

$q1.Range("A19").Value = 17
$q1.Range("B19").NumberFormat = "@"
$q1.Range("B19").Value = "005117"
$q1.Range("C19").Value = "金信价值精选灵活配置混合A"
$q1.Range("D19").NumberFormat = "@"
$q1.Range("D19").Value = "0.83"
$q1.Range("E19").NumberFormat = "@"
$q1.Range("E19").Value = "86.42"
$q1.Range("F19").NumberFormat = "@"
$q1.Range("F19").Value = "2.75"
$q1.Range("G19").NumberFormat = "@"
$q1.Range("G19").Value = "0.0228"
$q1.Range("H19").Value = 9

$q1.Range("A20").Value = 18
$q1.Range("B20").NumberFormat = "@"
$q1.Range("B20").Value = "519656"
$q1.Range("C20").Value = "银河灵活配置混合 - A"
$q1.Range("D20").NumberFormat = "@"
$q1.Range("D20").Value = "0.72"
$q1.Range("E20").NumberFormat = "@"
$q1.Range("E20").Value = "59.27"
$q1.Range("F20").NumberFormat = "@"
$q1.Range("F20").Value = "3.14"
$q1.Range("G20").NumberFormat = "@"
$q1.Range("G20").Value = "0.0226"
$q1.Range("H20").Value = 7

$q1.Range("A21").Value = 19
$q1.Range("B21").NumberFormat = "@"
$q1.Range("B21").Value = "010151"
$q1.Range("C21").Value = "西藏东财消费精选混合A"
$q1.Range("D21").NumberFormat = "@"
$q1.Range("D21").Value = "0.80"
$q1.Range("E21").NumberFormat = "@"
$q1.Range("E21").Value = "77.66"
$q1.Range("F21").NumberFormat = "@"
$q1.Range("F21").Value = "2.28"
$q1.Range("G21").NumberFormat = "@"
$q1.Range("G21").Value = "0.0182"
$q1.Range("H21").Value = 9

$q1.Range("A22").Value = 20
$q1.Range("B22").NumberFormat = "@"
$q1.Range("B22").Value = "002149"
$q1.Range("C22").Value = "嘉实新优选灵活配置混合"
$q1.Range("D22").NumberFormat = "@"
$q1.Range("D22").Value = "0.22"
$q1.Range("E22").NumberFormat = "@"
$q1.Range("E22").Value = "93.76"
$q1.Range("F22").NumberFormat = "@"
$q1.Range("F22").Value = "5.96"
$q1.Range("G22").NumberFormat = "@"
$q1.Range("G22").Value = "0.0131"
$q1.Range("H22").Value = 6

$q1.Range("A23").Value = 21
$q1.Range("B23").NumberFormat = "@"
$q1.Range("B23").Value = "004301"
$q1.Range("C23").Value = "国寿安保稳信混合A"
$q1.Range("D23").NumberFormat = "@"
$q1.Range("D23").Value = "1.50"
$q1.Range("E23").NumberFormat = "@"
$q1.Range("E23").Value = "20.03"
$q1.Range("F23").NumberFormat = "@"
$q1.Range("F23").Value = "0.75"
$q1.Range("G23").NumberFormat = "@"
$q1.Range("G23").Value = "0.0112"
$q1.Range("H23").Value = 10

$q1.Range("A24").Value = 22
$q1.Range("B24").NumberFormat = "@"
$q1.Range("B24").Value = "005572"
$q1.Range("C24").Value = "中银证券新能源灵活配置混合C"
$q1.Range("D24").NumberFormat = "@"
$q1.Range("D24").Value = "0.28"
$q1.Range("E24").NumberFormat = "@"
$q1.Range("E24").Value = "90.25"
$q1.Range("F24").NumberFormat = "@"
$q1.Range("F24").Value = "3.97"
$q1.Range("G24").NumberFormat = "@"
$q1.Range("G24").Value = "0.0111"
$q1.Range("H24").Value = 8

$q1.Range("A25").Value = 23
$q1.Range("B25").NumberFormat = "@"
$q1.Range("B25").Value = "002810"
$q1.Range("C25").Value = "金信转型创新成长灵活配置混合"
$q1.Range("D25").NumberFormat = "@"
$q1.Range("D25").Value = "0.18"
$q1.Range("E25").NumberFormat = "@"
$q1.Range("E25").Value = "81.12"
$q1.Range("F25").NumberFormat = "@"
$q1.Range("F25").Value = "5.75"
$q1.Range("G25").NumberFormat = "@"
$q1.Range("G25").Value = "0.0104"
$q1.Range("H25").Value = 2

$q1.Range("A26").Value = 24
$q1.Range("B26").NumberFormat = "@"
$q1.Range("B26").Value = "519657"
$q1.Range("C26").Value = "银河灵活配置混合 - C"
$q1.Range("D26").NumberFormat = "@"
$q1.Range("D26").Value = "0.33"
$q1.Range("E26").NumberFormat = "@"
$q1.Range("E26").Value = "59.27"
$q1.Range("F26").NumberFormat = "@"
$q1.Range("F26").Value = "3.14"
$q1.Range("G26").NumberFormat = "@"
$q1.Range("G26").Value = "0.0104"
$q1.Range("H26").Value = 7

$q1.Range("A27").Value = 25
$q1.Range("B27").NumberFormat = "@"
$q1.Range("B27").Value = "007578"
$q1.Range("C27").Value = "宝盈新锐灵活配置混合C"
$q1.Range("D27").NumberFormat = "@"
$q1.Range("D27").Value = "0.20"
$q1.Range("E27").NumberFormat = "@"
$q1.Range("E27").Value = "93.26"
$q1.Range("F27").NumberFormat = "@"
$q1.Range("F27").Value = "4.89"
$q1.Range("G27").NumberFormat = "@"
$q1.Range("G27").Value = "0.0098"
$q1.Range("H27").Value = 5

$q1.Range("A28").Value = 26
$q1.Range("B28").NumberFormat = "@"
$q1.Range("B28").Value = "010152"
$q1.Range("C28").Value = "西藏东财消费精选混合C"
$q1.Range("D28").NumberFormat = "@"
$q1.Range("D28").Value = "0.32"
$q1.Range("E28").NumberFormat = "@"
$q1.Range("E28").Value = "77.66"
$q1.Range("F28").NumberFormat = "@"
$q1.Range("F28").Value = "2.28"
$q1.Range("G28").NumberFormat = "@"
$q1.Range("G28").Value = "0.0073"
$q1.Range("H28").Value = 9

$q1.Range("A29").Value = 27
$q1.Range("B29").NumberFormat = "@"
$q1.Range("B29").Value = "519617"
$q1.Range("C29").Value = "银河君信灵活配置混合C"
$q1.Range("D29").NumberFormat = "@"
$q1.Range("D29").Value = "0.64"
$q1.Range("E29").NumberFormat = "@"
$q1.Range("E29").Value = "24.42"
$q1.Range("F29").NumberFormat = "@"
$q1.Range("F29").Value = "0.74"
$q1.Range("G29").NumberFormat = "@"
$q1.Range("G29").Value = "0.0047"
$q1.Range("H29").Value = 8

$q1.Range("A30").Value = 28
$q1.Range("B30").NumberFormat = "@"
$q1.Range("B30").Value = "005118"
$q1.Range("C30").Value = "金信价值精选灵活配置混合C"
$q1.Range("D30").NumberFormat = "@"
$q1.Range("D30").Value = "0.05"
$q1.Range("E30").NumberFormat = "@"
$q1.Range("E30").Value = "86.42"
$q1.Range("F30").NumberFormat = "@"
$q1.Range("F30").Value = "2.75"
$q1.Range("G30").NumberFormat = "@"
$q1.Range("G30").Value = "0.0014"
$q1.Range("H30").Value = 9

$q1.Range("A31").Value = 29
$q1.Range("B31").NumberFormat = "@"
$q1.Range("B31").Value = "004302"
$q1.Range("C31").Value = "国寿安保稳信混合C"
$q1.Range("D31").NumberFormat = "@"
$q1.Range("D31").Value = "0.01"
$q1.Range("E31").NumberFormat = "@"
$q1.Range("E31").Value = "20.03"
$q1.Range("F31").NumberFormat = "@"
$q1.Range("F31").Value = "0.75"
$q1.Range("G31").NumberFormat = "@"
$q1.Range("G31").Value = "0.0001"
$q1.Range("H31").Value = 10

$q1.Range("A32").Value = 30
$q1.Range("B32").NumberFormat = "@"
$q1.Range("B32").Value = "015406"
$q1.Range("C32").Value = "国寿安保稳信混合E"
$q1.Range("D32").NumberFormat = "@"
$q1.Range("D32").Value = "0.00"
$q1.Range("E32").NumberFormat = "@"
$q1.Range("E32").Value = "20.03"
$q1.Range("F32").NumberFormat = "@"
$q1.Range("F32").Value = "0.75"
$q1.Range("G32").Value = 0
$q1.Range("H32").Value = 10

# 2) Prepend the new 2022-Q1 summary row to the "总计" sheet, shifting
#    the existing quarters down by one row
$ws = $wb.Worksheets.Item("总计")
$ws.Rows.Item(2).Insert()

# Copy the A-column index cell formatting down into the newly
# inserted row, then clear the stray formatting Insert() left on B:D
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("B2:D2").ClearFormats()

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q1"
$ws.Range("C2").Value = 31
$ws.Range("D2").Value = 4.44

# Renumber the index column for the rows that shifted down
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
